$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 11
$ws.Range("B5").Value = 7
$ws.Range("B6").Value = 4
$ws.Range("B12").Value = -8
$ws.Range("B25").Value = 9
$ws.Range("B27").Value = 4
$ws.Range("B30").Value = 5
